# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (Damasco - Dina) at the top of this variety's
# block, pushing the existing O'Higgins / Modesto / Dina rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 98; existing rows 98-104 shift down to 101-107.
$ws.Rows("98:100").Insert()

# New row 98: Damasco, Dina, Especial
$ws.Cells.Item(98, 1).Value2 = 8
$ws.Cells.Item(98, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(98, 3).Value2 = "Coquimbo"
$ws.Cells.Item(98, 4).Value2 = 44931
$ws.Cells.Item(98, 5).Value2 = 4
$ws.Cells.Item(98, 6).Value2 = "Fruta"
$ws.Cells.Item(98, 7).Value2 = 100103
$ws.Cells.Item(98, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(98, 9).Value2 = 100103003
$ws.Cells.Item(98, 10).Value2 = "Damasco"
$ws.Cells.Item(98, 11).Value2 = "Dina"
$ws.Cells.Item(98, 12).Value2 = "Especial"
$ws.Cells.Item(98, 13).Value2 = 300
$ws.Cells.Item(98, 14).Value2 = 22000
$ws.Cells.Item(98, 15).Value2 = 23000
$ws.Cells.Item(98, 16).Value2 = 22500
$ws.Cells.Item(98, 17).Value2 = "$/caja 16 kilos"
$ws.Cells.Item(98, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(98, 19).Value2 = 1406
$ws.Cells.Item(98, 20).Value2 = 16

# New row 99: Damasco, Dina, Primera
$ws.Cells.Item(99, 1).Value2 = 8
$ws.Cells.Item(99, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(99, 3).Value2 = "Coquimbo"
$ws.Cells.Item(99, 4).Value2 = 44931
$ws.Cells.Item(99, 5).Value2 = 4
$ws.Cells.Item(99, 6).Value2 = "Fruta"
$ws.Cells.Item(99, 7).Value2 = 100103
$ws.Cells.Item(99, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(99, 9).Value2 = 100103003
$ws.Cells.Item(99, 10).Value2 = "Damasco"
$ws.Cells.Item(99, 11).Value2 = "Dina"
$ws.Cells.Item(99, 12).Value2 = "Primera"
$ws.Cells.Item(99, 13).Value2 = 200
$ws.Cells.Item(99, 14).Value2 = 19000
$ws.Cells.Item(99, 15).Value2 = 20000
$ws.Cells.Item(99, 16).Value2 = 19500
$ws.Cells.Item(99, 17).Value2 = "$/caja 16 kilos"
$ws.Cells.Item(99, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(99, 19).Value2 = 1219
$ws.Cells.Item(99, 20).Value2 = 16

# New row 100: Damasco, Dina, Segunda
$ws.Cells.Item(100, 1).Value2 = 8
$ws.Cells.Item(100, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(100, 3).Value2 = "Coquimbo"
$ws.Cells.Item(100, 4).Value2 = 44931
$ws.Cells.Item(100, 5).Value2 = 4
$ws.Cells.Item(100, 6).Value2 = "Fruta"
$ws.Cells.Item(100, 7).Value2 = 100103
$ws.Cells.Item(100, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(100, 9).Value2 = 100103003
$ws.Cells.Item(100, 10).Value2 = "Damasco"
$ws.Cells.Item(100, 11).Value2 = "Dina"
$ws.Cells.Item(100, 12).Value2 = "Segunda"
$ws.Cells.Item(100, 13).Value2 = 300
$ws.Cells.Item(100, 14).Value2 = 15000
$ws.Cells.Item(100, 15).Value2 = 16000
$ws.Cells.Item(100, 16).Value2 = 15500
$ws.Cells.Item(100, 17).Value2 = "$/caja 16 kilos"
$ws.Cells.Item(100, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(100, 19).Value2 = 969
$ws.Cells.Item(100, 20).Value2 = 16
